$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-average-scheduled-hours"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# The "Fixed Value" cell for Extension.url shares the same URL string as
# the Metadata URL cell, so it must be updated to match too.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-average-scheduled-hours"

# Clear the Constraint(s) value on the base "Extension" row (row 2); the
# ele-1/ext-1 invariant text now only lives on the Extension.extension row.
$elements.Range("AI2").Value = ""
